$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.164309024810791
$ws.Range("B1").Value = 2.421557903289795
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.376310110092163
$ws.Range("E1").Value = 1.235073208808899
